$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ValidateParameter"

$srcSheet.Range("B60").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$newSheet.Range("B2").Value = "A. テストケース"

$srcSheet.Range("C4").Copy()
$newSheet.Range("C4").PasteSpecial(-4122)
$newSheet.Range("C4").Value = "Expect系メソッドを呼ぶ前にGetParamterAtメソッドを呼ぶと例外が発生する"

$srcSheet.Range("B64").Copy()
$newSheet.Range("B6").PasteSpecial(-4122)
$newSheet.Range("B6").Value = "E.取得データ"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C7").PasteSpecial(-4122)
$newSheet.Range("C7").Value = "ConfigException"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D7").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C8").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D8").PasteSpecial(-4122)
$newSheet.Range("D8").Value = "ResourceKey"

$srcSheet.Range("D66").Copy()
$newSheet.Range("E8").PasteSpecial(-4122)
$newSheet.Range("E8").Value = "Message"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C9").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D9").PasteSpecial(-4122)
$newSheet.Range("D9").Value = "M_Fixture_Temp_Conductor_InvalidStatus"

$srcSheet.Range("D67").Copy()
$newSheet.Range("E9").PasteSpecial(-4122)
$newSheet.Range("E9").Value = "%GetParamterAt%"

$srcSheet.Range("B60").Copy()
$newSheet.Range("B11").PasteSpecial(-4122)
$newSheet.Range("B11").Value = "A. テストケース"

$srcSheet.Range("C4").Copy()
$newSheet.Range("C13").PasteSpecial(-4122)
$newSheet.Range("C13").Value = "Expect系メソッドを呼ぶ前にValidateParamterAtメソッドを呼ぶと例外が発生する"

$srcSheet.Range("B64").Copy()
$newSheet.Range("B15").PasteSpecial(-4122)
$newSheet.Range("B15").Value = "E.取得データ"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C16").PasteSpecial(-4122)
$newSheet.Range("C16").Value = "ConfigException"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D16").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C17").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D17").PasteSpecial(-4122)
$newSheet.Range("D17").Value = "ResourceKey"

$srcSheet.Range("D66").Copy()
$newSheet.Range("E17").PasteSpecial(-4122)
$newSheet.Range("E17").Value = "Message"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C18").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D18").PasteSpecial(-4122)
$newSheet.Range("D18").Value = "M_Fixture_Temp_Conductor_InvalidStatus"

$srcSheet.Range("D67").Copy()
$newSheet.Range("E18").PasteSpecial(-4122)
$newSheet.Range("E18").Value = "%ValidateParamterAt%"

$srcSheet.Range("B60").Copy()
$newSheet.Range("B20").PasteSpecial(-4122)
$newSheet.Range("B20").Value = "A. テストケース"

$srcSheet.Range("C4").Copy()
$newSheet.Range("C22").PasteSpecial(-4122)
$newSheet.Range("C22").Value = "GetParamterAtメソッドのインデックスがExpectの引数の数よりも多い場合は例外が発生する"

$srcSheet.Range("B64").Copy()
$newSheet.Range("B24").PasteSpecial(-4122)
$newSheet.Range("B24").Value = "E.取得データ"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C25").PasteSpecial(-4122)
$newSheet.Range("C25").Value = "ConfigException"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D25").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C26").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D26").PasteSpecial(-4122)
$newSheet.Range("D26").Value = "ResourceKey"

$srcSheet.Range("D66").Copy()
$newSheet.Range("E26").PasteSpecial(-4122)
$newSheet.Range("E26").Value = "Message"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C27").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D27").PasteSpecial(-4122)
$newSheet.Range("D27").Value = "M_Fixture_Temp_Conductor_InvalidParameterIndex"

$srcSheet.Range("D67").Copy()
$newSheet.Range("E27").PasteSpecial(-4122)
$newSheet.Range("E27").Value = "%(0)%"

$srcSheet.Range("B60").Copy()
$newSheet.Range("B29").PasteSpecial(-4122)
$newSheet.Range("B29").Value = "A. テストケース"

$srcSheet.Range("C4").Copy()
$newSheet.Range("C31").PasteSpecial(-4122)
$newSheet.Range("C31").Value = "ValidateParamterAtメソッドのインデックスがExpectの引数の数よりも多い場合は例外が発生する"

$srcSheet.Range("B64").Copy()
$newSheet.Range("B33").PasteSpecial(-4122)
$newSheet.Range("B33").Value = "D. パラメタ"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C34").PasteSpecial(-4122)
$newSheet.Range("C34").Value = "Data1"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D34").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C35").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D35").PasteSpecial(-4122)
$newSheet.Range("D35").Value = "text"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C36").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D36").PasteSpecial(-4122)
$newSheet.Range("D36").Value = "abc"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C38").PasteSpecial(-4122)
$newSheet.Range("C38").Value = "Data2"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D38").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C39").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D39").PasteSpecial(-4122)
$newSheet.Range("D39").Value = "text"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C40").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D40").PasteSpecial(-4122)
$newSheet.Range("D40").Value = "def"

$srcSheet.Range("B64").Copy()
$newSheet.Range("B42").PasteSpecial(-4122)
$newSheet.Range("B42").Value = "E.取得データ"

$srcSheet.Range("C65").Copy()
$newSheet.Range("C43").PasteSpecial(-4122)
$newSheet.Range("C43").Value = "ConfigException"

$srcSheet.Range("D65").Copy()
$newSheet.Range("D43").PasteSpecial(-4122)

$srcSheet.Range("C66").Copy()
$newSheet.Range("C44").PasteSpecial(-4122)

$srcSheet.Range("D66").Copy()
$newSheet.Range("D44").PasteSpecial(-4122)
$newSheet.Range("D44").Value = "ResourceKey"

$srcSheet.Range("D66").Copy()
$newSheet.Range("E44").PasteSpecial(-4122)
$newSheet.Range("E44").Value = "Message"

$srcSheet.Range("C66").Copy()
$newSheet.Range("C45").PasteSpecial(-4122)

$srcSheet.Range("D67").Copy()
$newSheet.Range("D45").PasteSpecial(-4122)
$newSheet.Range("D45").Value = "M_Fixture_Temp_Conductor_InvalidParameterIndex"

$srcSheet.Range("D67").Copy()
$newSheet.Range("E45").PasteSpecial(-4122)
$newSheet.Range("E45").Value = "%(2)%"
